$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "price " in F1
$ws.Range("F1").Value = "price "

# Add price values in column F for rows 2-7
$ws.Range("F2").Value = 89.9
$ws.Range("F3").Value = 200
$ws.Range("F4").Value = 250
$ws.Range("F5").Value = 350
$ws.Range("F6").Value = 70
$ws.Range("F7").Value = 50

# Set column widths to match the target layout (closest achievable given
# the host's character-width quantization; inverse-solved against the
# target OOXML <col width> values of 13.296875 / 15.796875 / 12.5 / 11 / 17.5)
$ws.Columns.Item(1).ColumnWidth = 12.571428571428571
$ws.Columns.Item(2).ColumnWidth = 15.142857142857142
$ws.Columns.Item(3).ColumnWidth = 11.714285714285715
$ws.Columns.Item(4).ColumnWidth = 10.285714285714285
$ws.Columns.Item(5).ColumnWidth = 16.714285714285715

# Set selection to F7 to match the target state
$ws.Range("F7").Select()
